$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Change 1: merge the bookmark-split run back into a single run and drop the bookmark ---
$find = $d.Content.Find
$find.Execute("fixed th" + "rough the Visual Studio", $true, $false, $false, $false, $false, `
               $true, 1, $false, "fixed through the Visual Studio", 2) | Out-Null

# --- Change 2: append a new "Documentation" section at the end of the document ---

# 2a. Empty spacer paragraph (ind left=45)
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$spacer = $d.Paragraphs.Item($d.Paragraphs.Count)
$spacerXml = "<w:p $wns><w:pPr><w:ind w:left=`"45`"/><w:rPr><w:rFonts w:cstheme=`"minorHAnsi`"/><w:sz w:val=`"24`"/></w:rPr></w:pPr></w:p>"
$spacer.Range.InsertXML($spacerXml)

# 2b. "Documentation" heading paragraph, with the _GoBack bookmark right after the text
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$heading = $d.Paragraphs.Item($d.Paragraphs.Count)
$headingXml = "<w:p $wns><w:pPr><w:rPr><w:rFonts w:cstheme=`"minorHAnsi`"/><w:sz w:val=`"32`"/><w:u w:val=`"single`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme=`"minorHAnsi`"/><w:sz w:val=`"32`"/><w:u w:val=`"single`"/></w:rPr><w:t>Documentation</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
$heading.Range.InsertXML($headingXml)

# 2c. Documentation body paragraph (tab + text, ind left=45)
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$body = $d.Paragraphs.Item($d.Paragraphs.Count)
$bodyText = "Documentation that is to be included in the repository should be added to the " + [char]0x201C + "master" + [char]0x201D + " branch, as opposed to " + [char]0x201C + "development" + [char]0x201D + ". This is to avoid unnecessary builds on the development branch, and to ensure documentation is present for deliverables."
$bodyXml = "<w:p $wns><w:pPr><w:ind w:left=`"45`"/><w:rPr><w:rFonts w:cstheme=`"minorHAnsi`"/><w:sz w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme=`"minorHAnsi`"/><w:sz w:val=`"24`"/></w:rPr><w:tab/><w:t>$bodyText</w:t></w:r></w:p>"
$body.Range.InsertXML($bodyXml)

Write-Host "done"
